# Add "keywords" and "url" columns to the "Tool" sheet schema, inserted
# right before the existing "id" column (which currently lives in C1).
#
# Before: A=developer_team, B=technical_area, C=id, D=name, E=description
# After:  A=developer_team, B=technical_area, C=keywords, D=url, E=id, F=name, G=description

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tool")

# Shift existing C:E columns (id, name, description) two columns to the right,
# making room at C1:D1 for the new headers.
$xlShiftToRight = -4161
$ws.Range("C1:D1").EntireColumn.Insert($xlShiftToRight)

$ws.Range("C1").Value = "keywords"
$ws.Range("D1").Value = "url"
